$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4485
$ws.Range("I62").Value = 4529.6665
$ws.Range("K62").Value = 4529.6665
$ws.Range("M62").Value = -3905.6665
$ws.Range("H65").Value = 4485
$ws.Range("I65").Value = 4529.6665
$ws.Range("K65").Value = 22648.3325
$ws.Range("M65").Value = -19528.3325
$ws.Range("H70").Value = 4487.909
$ws.Range("I70").Value = 3625
$ws.Range("J70").Value = 4574.2
$ws.Range("K70").Value = 10875
$ws.Range("L70").Value = 13722.6
$ws.Range("M70").Value = -10605
$ws.Range("N70").Value = -14262.6
$ws.Range("H73").Value = 4487.909
$ws.Range("I73").Value = 3625
$ws.Range("J73").Value = 4574.2
$ws.Range("K73").Value = 10875
$ws.Range("L73").Value = 13722.6
$ws.Range("M73").Value = -9939
$ws.Range("N73").Value = -15594.6
$ws.Range("H80").Value = 1240
$ws.Range("I80").Value = 1300
$ws.Range("J80").Value = 1000
$ws.Range("K80").Value = 3900
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -2902
$ws.Range("N80").Value = -4996
$ws.Range("H83").Value = 1240
$ws.Range("I83").Value = 1300
$ws.Range("J83").Value = 1000
$ws.Range("K83").Value = 11700
$ws.Range("L83").Value = 9000
$ws.Range("M83").Value = -6708
$ws.Range("N83").Value = -18984
$ws.Range("H86").Value = 4232.273
$ws.Range("I86").Value = 4371.75
$ws.Range("J86").Value = 3860.3333
$ws.Range("K86").Value = 4371.75
$ws.Range("L86").Value = 3860.3333
$ws.Range("M86").Value = -3248.75
$ws.Range("N86").Value = -6106.3333
$ws.Range("H88").Value = 810.5
$ws.Range("I88").Value = 360
$ws.Range("J88").Value = 1080.8
$ws.Range("K88").Value = 360
$ws.Range("L88").Value = 1080.8
$ws.Range("M88").Value = 46
$ws.Range("N88").Value = -1892.8
$ws.Range("H89").Value = 4232.273
$ws.Range("I89").Value = 4371.75
$ws.Range("J89").Value = 3860.3333
$ws.Range("K89").Value = 21858.75
$ws.Range("L89").Value = 19301.6665
$ws.Range("M89").Value = -16242.75
$ws.Range("N89").Value = -30533.6665
$ws.Range("H91").Value = 810.5
$ws.Range("I91").Value = 360
$ws.Range("J91").Value = 1080.8
$ws.Range("K91").Value = 360
$ws.Range("L91").Value = 1080.8
$ws.Range("M91").Value = 1044
$ws.Range("N91").Value = -3888.8
$ws.Range("H98").Value = 2364.762
$ws.Range("J98").Value = 3227.8333
$ws.Range("L98").Value = 3227.8333
$ws.Range("N98").Value = -6223.8333
$ws.Range("H106").Value = 2545.923
$ws.Range("I106").Value = 2417.0908
$ws.Range("J106").Value = 3254.5
$ws.Range("K106").Value = 2417.0908
$ws.Range("L106").Value = 3254.5
$ws.Range("M106").Value = -1786.0908
$ws.Range("N106").Value = -4516.5
$ws.Range("H112").Value = 2221.077
$ws.Range("J112").Value = 2264.5833
$ws.Range("L112").Value = 6793.749899999999
$ws.Range("N112").Value = -9009.749899999999
$ws.Range("H122").Value = 2364.762
$ws.Range("J122").Value = 3227.8333
$ws.Range("L122").Value = 9683.499899999999
$ws.Range("N122").Value = -14583.4999
$ws.Range("H129").Value = 3068.7058
$ws.Range("I129").Value = 1137.8
$ws.Range("J129").Value = 3873.25
$ws.Range("K129").Value = 3413.4
$ws.Range("L129").Value = 11619.75
$ws.Range("M129").Value = 1586.6
$ws.Range("N129").Value = -21619.75
$ws.Range("H132").Value = 1999.7059
$ws.Range("I132").Value = 1941.5
$ws.Range("K132").Value = 5824.5
$ws.Range("M132").Value = -3294.5

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2322.7
$ws.Range("I2").Value = 2404
$ws.Range("K2").Value = 2404
$ws.Range("M2").Value = -2291
$ws.Range("H63").Value = 2530.5557
$ws.Range("I63").Value = 1692.25
$ws.Range("J63").Value = 3201.2
$ws.Range("K63").Value = 1692.25
$ws.Range("L63").Value = 3201.2
$ws.Range("M63").Value = -1006.25
$ws.Range("N63").Value = -4573.2
$ws.Range("H66").Value = 2530.5557
$ws.Range("I66").Value = 1692.25
$ws.Range("J66").Value = 3201.2
$ws.Range("K66").Value = 8461.25
$ws.Range("L66").Value = 16006
$ws.Range("M66").Value = -5029.25
$ws.Range("N66").Value = -22870
$ws.Range("H88").Value = 2317.2222
$ws.Range("I88").Value = 753.2
$ws.Range("J88").Value = 2918.7693
$ws.Range("K88").Value = 753.2
$ws.Range("L88").Value = 2918.7693
$ws.Range("M88").Value = -347.2
$ws.Range("N88").Value = -3730.7693
$ws.Range("H91").Value = 2317.2222
$ws.Range("I91").Value = 753.2
$ws.Range("J91").Value = 2918.7693
$ws.Range("K91").Value = 753.2
$ws.Range("L91").Value = 2918.7693
$ws.Range("M91").Value = 650.8
$ws.Range("N91").Value = -5726.7693
$ws.Range("H116").Value = 2322.7
$ws.Range("I116").Value = 2404
$ws.Range("K116").Value = 2404
$ws.Range("M116").Value = -110

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2322.7
$ws.Range("I3").Value = 2404
$ws.Range("K3").Value = 2404
$ws.Range("M3").Value = -2290
$ws.Range("H20").Value = 1321.5714
$ws.Range("I20").Value = 692.75
$ws.Range("K20").Value = 692.75
$ws.Range("M20").Value = -445.75

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2199.15
$ws.Range("I58").Value = 2181.3125
$ws.Range("J58").Value = 2270.5
$ws.Range("K58").Value = 2181.3125
$ws.Range("L58").Value = 2270.5
$ws.Range("M58").Value = -1978.3125
$ws.Range("N58").Value = -2676.5
$ws.Range("H132").Value = 1596.2858
$ws.Range("I132").Value = 1435
$ws.Range("K132").Value = 4305
$ws.Range("M132").Value = -1775
$ws.Range("H136").Value = 2199.15
$ws.Range("I136").Value = 2181.3125
$ws.Range("J136").Value = 2270.5
$ws.Range("K136").Value = 6543.9375
$ws.Range("L136").Value = 6811.5
$ws.Range("M136").Value = -3993.9375
$ws.Range("N136").Value = -11911.5

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2009.8148
$ws.Range("I4").Value = 1910.8823
$ws.Range("J4").Value = 2178
$ws.Range("K4").Value = 5732.6469
$ws.Range("L4").Value = 6534
$ws.Range("M4").Value = -5620.6469
$ws.Range("N4").Value = -6758
$ws.Range("H97").Value = 825.25
$ws.Range("I97").Value = 1333.8
$ws.Range("J97").Value = 462
$ws.Range("K97").Value = 4001.4
$ws.Range("L97").Value = 1386
$ws.Range("M97").Value = -3505.4
$ws.Range("N97").Value = -2378
$ws.Range("H98").Value = 126.666664
$ws.Range("I98").Value = 80
$ws.Range("J98").Value = 150
$ws.Range("K98").Value = 240
$ws.Range("L98").Value = 450
$ws.Range("M98").Value = 1258
$ws.Range("N98").Value = -3446
$ws.Range("H113").Value = 658
$ws.Range("I113").Value = 577.6
$ws.Range("J113").Value = 698.2
$ws.Range("K113").Value = 1732.8
$ws.Range("L113").Value = 2094.6
$ws.Range("M113").Value = 437.1999999999998
$ws.Range("N113").Value = -6434.6

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 10421968
$ws.Range("I122").Value = 10421968
$ws.Range("K122").Value = 31265904
$ws.Range("M122").Value = -31263454
$ws.Range("H126").Value = 3166.6667
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").Value = ""
$ws.Range("H132").Value = 1865.8334
$ws.Range("I132").Value = 1839
$ws.Range("K132").Value = 5517
$ws.Range("M132").Value = -2987

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3501.3333
$ws.Range("I122").Value = 3501.3333
$ws.Range("K122").Value = 10503.9999
$ws.Range("M122").Value = -8053.999899999999

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 25000
$ws.Range("J94").Value = 25000
$ws.Range("L94").Value = 25000
$ws.Range("N94").Value = -26802
$ws.Range("H110").Value = 160599.5
$ws.Range("J110").Value = 160599.5
$ws.Range("L110").Value = 160599.5
$ws.Range("N110").Value = -168779.5
